$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "585.98") are forced
# to Text format first. Excel normally auto-detects such strings and stores
# them as numeric values, but the source data keeps them as literal text
# (this also preserves exact formatting such as trailing zeros, e.g. "0.650").

# Row 2
$ws.Range("D2").Value = '67.206.37'
$ws.Range("E2").Value = '  +4.54%  '
# Row 3
$ws.Range("D3").Value = '3.474.23'
$ws.Range("E3").Value = '  +4.39%  '
# Row 4
$ws.Range("E4").Value = '  +0.02%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.98'
$ws.Range("E5").Value = '  +5.94%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.79'
$ws.Range("E6").Value = '  +7.87%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.632'
$ws.Range("E7").Value = '  +0.48%  '
# Row 8
$ws.Range("D8").Value = '3.471.71'
$ws.Range("E8").Value = '  +4.55%  '
# Row 10
$ws.Range("E10").Value = '  -0.08%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.650'
$ws.Range("E11").Value = '  +2.58%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.31'
$ws.Range("E12").Value = '  +6.11%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000279'
$ws.Range("E13").Value = '  -0.03%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.44'
$ws.Range("E14").Value = '  +4.03%  '
# Row 15
$ws.Range("D15").Value = '4.037.48'
$ws.Range("E15").Value = '  +4.57%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.83'
$ws.Range("E16").Value = '  +3.96%  '
# Row 17
$ws.Range("D17").Value = '3.473.95'
$ws.Range("E17").Value = '  +4.32%  '
# Row 18
$ws.Range("D18").Value = '67.154.53'
$ws.Range("E18").Value = '  +4.35%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.18'
$ws.Range("E19").Value = '  +4.02%  '
# Row 20
$ws.Range("E20").Value = '  -1.40%  '
# Row 21
$ws.Range("E21").Value = '  +3.54%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '486.52'
$ws.Range("E22").Value = '  +7.62%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.38'
$ws.Range("E23").Value = '  +7.40%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.81'
$ws.Range("E24").Value = '  +21.07%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.49'
$ws.Range("E25").Value = '  +10.64%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.85'
$ws.Range("E26").Value = '  +2.48%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.96'
$ws.Range("E27").Value = '  +3.20%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.97'
$ws.Range("E28").Value = '  +4.21%  '
# Row 29
$ws.Range("E29").Value = '  +6.62%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.42'
$ws.Range("E30").Value = '  +1.22%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.22'
$ws.Range("E31").Value = '  +11.07%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '602.96'
$ws.Range("E32").Value = '  +5.44%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.77'
$ws.Range("E33").Value = '  +3.28%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.83'
$ws.Range("E34").Value = '  +2.52%  '
# Row 35
$ws.Range("E35").Value = '  +4.93%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.148'
$ws.Range("E36").Value = '  +4.57%  '
# Row 37
$ws.Range("E37").Value = '  -0.04%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.67'
$ws.Range("E38").Value = '  +3.73%  '
# Row 39
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.387'
$ws.Range("E39").Value = '  +5.61%  '
# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.54'
$ws.Range("E40").Value = '  +0.59%  '
# Row 41
$ws.Range("D41").Value = '0.0₃0761'
$ws.Range("E41").Value = '  +4.43%  '
# Row 42
$ws.Range("D42").Value = '3.267.29'
$ws.Range("E42").Value = '  +6.36%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.91'
$ws.Range("E43").Value = '  +6.24%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0430'
$ws.Range("E44").Value = '  +4.06%  '
# Row 45
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  +24.97%  '
# Row 46
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.54'
$ws.Range("E46").Value = '  +3.67%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.25'
$ws.Range("E47").Value = '  +1.89%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.135'
$ws.Range("E48").Value = '  +1.31%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.26'
$ws.Range("E49").Value = '  +11.45%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.76'
$ws.Range("E50").Value = '  +7.49%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.11%  '
